# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rows 16-18 on "Hoja1" hold a small table of worker account-statement
# rows (Tipo Doc / N Doc / Nombre / Periodo Mora / Valor Mora / Salario
# Basico). The data for "ANUARD JOEL SALAS CARDENAS" (period 2410) and
# "MARIA LORENA RODRIGUEZ CARABALLO" (period 2409) is re-ordered and the
# Salario Basico for Maria Lorena is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Anuard, period 2409 (was 2410), Valor Mora 26000 (was 36400)
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1002243950"
$ws.Cells.Item(16, 4).Value = "ANUARD JOEL SALAS CARDENAS"
$ws.Cells.Item(16, 5).Value = "2409"
$ws.Cells.Item(16, 6).Value = 26000
$ws.Cells.Item(16, 7).Value = 1300000

# Row 17: Maria Lorena, period 2409, Valor Mora 26000, Salario Basico 1423500 (was 1300000)
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1001900040"
$ws.Cells.Item(17, 4).Value = "MARIA LORENA RODRIGUEZ CARABALLO"
$ws.Cells.Item(17, 5).Value = "2409"
$ws.Cells.Item(17, 6).Value = 26000
$ws.Cells.Item(17, 7).Value = 1423500

# Row 18: Anuard, period 2410 (was Maria Lorena / 2409), Valor Mora 36400 (was 26000)
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1002243950"
$ws.Cells.Item(18, 4).Value = "ANUARD JOEL SALAS CARDENAS"
$ws.Cells.Item(18, 5).Value = "2410"
$ws.Cells.Item(18, 6).Value = 36400
$ws.Cells.Item(18, 7).Value = 1300000
